# Append the new run-log row (row 12) to Sheet1, mirroring the formatting
# of the preceding rows (row 11 in particular).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 12

# Copy the style/format of row 11 down into row 12 first, so the new row
# matches the look (borders/alignment/number formats) of the existing log
# rows.
$ws.Range("A11:H11").Copy()
$ws.Range("A12:H12").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Cells.Item($newRow, 1).Value = "2025-08-14 09:43:02 UTC"
$ws.Cells.Item($newRow, 2).Value = "2025-08-14 15:13:02 IST"
$ws.Cells.Item($newRow, 3).Value = "SKIPPED"
$ws.Cells.Item($newRow, 4).Value = "No change in PDF. Skipping download & Excel update."
$ws.Cells.Item($newRow, 5).Value = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-08-2025.pdf"
$ws.Cells.Item($newRow, 6).Value = ""
$ws.Cells.Item($newRow, 7).Value = 0
$ws.Cells.Item($newRow, 8).Value = ""
